$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.629.73"
$ws.Range("E2").Value = "  -1.64%  "

$ws.Range("D3").Value = "2.585.66"
$ws.Range("E3").Value = "  -2.52%  "

$ws.Range("E4").Value = "  +0.43%  "

$ws.Range("D5").Value = "'508.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.83%  "

$ws.Range("D6").Value = "'156.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.63%  "

$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("E8").Value = "  -5.75%  "

$ws.Range("D9").Value = "2.591.96"
$ws.Range("E9").Value = "  -3.40%  "

$ws.Range("D10").Value = "'6.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.87%  "

$ws.Range("D11").Value = "'0.104"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.10%  "

$ws.Range("D12").Value = "'0.347"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.28%  "

$ws.Range("E13").Value = "  +0.85%  "

$ws.Range("D14").Value = "3.034.01"
$ws.Range("E14").Value = "  -3.29%  "

$ws.Range("D15").Value = "60.625.61"
$ws.Range("E15").Value = "  -0.91%  "

$ws.Range("D16").Value = "'21.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.26%  "

$ws.Range("E17").Value = "  -2.21%  "

$ws.Range("D18").Value = "2.582.57"
$ws.Range("E18").Value = "  -3.58%  "

$ws.Range("D19").Value = "'4.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.75%  "

$ws.Range("D20").Value = "'347.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.72%  "

$ws.Range("D21").Value = "'10.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.67%  "

$ws.Range("D22").Value = "'6.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.44%  "

$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("D24").Value = "'60.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.25%  "

$ws.Range("D25").Value = "'0.421"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.65%  "

$ws.Range("E26").Value = "  -1.85%  "

$ws.Range("D27").Value = "2.695.83"
$ws.Range("E27").Value = "  -3.82%  "

$ws.Range("D28").Value = "'0.996"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.42%  "

$ws.Range("D29").Value = "0.0₃0850"
$ws.Range("E29").Value = "  -2.75%  "

$ws.Range("D30").Value = "'7.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.35%  "

$ws.Range("E31").Value = "  +0.15%  "

$ws.Range("D32").Value = "'19.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.34%  "

$ws.Range("D33").Value = "'153.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.69%  "

$ws.Range("D34").Value = "'1.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.60%  "

$ws.Range("D35").Value = "'5.74"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.22%  "

$ws.Range("D36").Value = "'4.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.74%  "

$ws.Range("D37").Value = "'1.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.62%  "

$ws.Range("D38").Value = "'0.855"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.62%  "

$ws.Range("D39").Value = "'1.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.63%  "

$ws.Range("D40").Value = "'0.849"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.20%  "

$ws.Range("D41").Value = "'36.25"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.43%  "

$ws.Range("D42").Value = "'3.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.42%  "

$ws.Range("D43").Value = "'298.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.37%  "

$ws.Range("D44").Value = "'0.623"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.64%  "

$ws.Range("D45").Value = "'0.100"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.74%  "

$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "'0.0560"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.96%  "

$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "'0.998"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.50%  "

$ws.Range("D48").Value = "'19.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.93%  "

$ws.Range("D49").Value = "'4.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.20%  "

$ws.Range("D50").Value = "'0.0234"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.89%  "

$ws.Range("D51").Value = "'10.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.29%  "

